$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds the last-changed date as an Excel serial date.
# All data rows (2 through 300) are being bumped from 45181 (2023-09-12)
# to 45182 (2023-09-13), i.e. the automated "last updated" stamp advanced
# by one day.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 300) { $lastRow = 300 }

$rng = $ws.Range("C2:C$lastRow")
$rng.Value = 45182
